# Update the "Chart" sheet notes at the bottom of the report template.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

# Bump the "last tested with" ReportServer version note (row 14, column A)
$ws.Cells.Item(14, 1).Value = "(Last tested with: ReportServer 4.1.0-6064) "

# Append a new note below it (row 15, column A) explaining the chart data ranges
$ws.Cells.Item(15, 1).Value = "(Note that in order to create the chart, we would need to know where the chart data is going appear to select the correct cell ranges for the axes.)"

# Reflect where the cursor/selection ended up after the edit
[void]$ws.Range("A12").Select()
